# Apply "Diferencia Stock" (column L) = "Stock Minimo Objetivo" (column K)
# for data rows 3 through 68, and update the summary cell C82
# (Total_Ajuste_Stock) to match C81 (Stock_Minimo_Objetivo).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

for ($r = 3; $r -le 68; $r++) {
    $kValue = $ws.Cells.Item($r, 11).Value2   # column K = 11
    $ws.Cells.Item($r, 12).Value2 = $kValue   # column L = 12
}

# Update Total_Ajuste_Stock (C82) to equal Stock_Minimo_Objetivo (C81)
$ws.Range("C82").Value2 = $ws.Range("C81").Value2
